$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K" header) values regenerated to use K instead of Strike# (rows 2-12)
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 4
$ws.Range("G10").Value = 7
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 0
